# [Receivable / Payable  Transactions] amount not show
#
# The "Amount" column on the invoice/transactions line template referenced the
# wrong merge field ({{currency amount}}) - it should reference the line's
# due_amount field ({{currency due_amount}}) so the amount actually renders.
#
# Also nudge the saved cursor position / column widths to match the
# as-saved state of the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core fix: Amount column template cell (F6) ---
# Was: {{currency amount}}  ->  Now: {{currency due_amount}}
$ws.Range("F6").Value = "{{currency due_amount}}"

# --- Cosmetic: selection/cursor position as saved in the authored file ---
[void]$ws.Range("F14").Select()

# --- Cosmetic: column widths drifted slightly on save; reproduce as closely
#     as this engine's character-width quantization allows. ---
$ws.Range("A1:C1").ColumnWidth = 21.6666
$ws.Columns.Item(4).ColumnWidth = 43.1667
$ws.Columns.Item(5).ColumnWidth = 24.6666
$ws.Columns.Item(6).ColumnWidth = 27.0
$ws.Columns.Item(7).ColumnWidth = 25.3333
$ws.Columns.Item(8).ColumnWidth = 7.5
